$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -11.95909999999999
$ws.Range("D6").Value = -7.908200000000001
$ws.Range("D7").Value = -7.548899999999992
$ws.Range("C8").Value = -12.25639999999999
$ws.Range("D8").Value = -7.901100000000002
$ws.Range("B12").Value = 5.675299999999998
$ws.Range("C12").Value = -14.77050000000003
$ws.Range("C14").Value = -12.463
$ws.Range("D19").Value = -8.342199999999993
$ws.Range("D21").Value = -7.669299999999998
$ws.Range("C22").Value = -11.1465
$ws.Range("D24").Value = -7.887299999999997
